# Apply the commit's changes to the workbook:
#  1. Insert a new "Player Info" sheet at the front with player bio data.
#  2. Rename the MATCH_CARD_LINK column to MATCH_CODE on the "ODI Batting"
#     and "ODI Bowling" sheets, and replace the full scorecard URLs with
#     just the bare match-code numbers.

$wb = $excel.ActiveWorkbook

# --- 1. Add the "Player Info" sheet as the new first sheet -----------------
# NB: grab the "ODI Batting" handle *after* Worksheets.Add() - fetching it
# beforehand and holding onto that variable across the Add() call ends up
# aliasing the newly inserted sheet in this host, not the original one.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo.Move($battingSheet)

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header look used by the other sheets.
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4664"
$playerInfo.Range("B2").Value = "Thangarasu Natarajan"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Medium"

# --- 2. Rename MATCH_CARD_LINK -> MATCH_CODE and trim the URLs down to the
#        bare match code on both "ODI Batting" and "ODI Bowling" sheets ----
# NB: this host's PowerShell binds positional args only - named args
# (-sheet ... -col ...) silently fail to bind, so call positionally.
function Update-MatchCodeColumn {
    param($sheet, $col, $lastRow)

    $headerCell = $sheet.Cells.Item(1, $col)
    $headerCell.Value = "MATCH_CODE"

    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $sheet.Cells.Item($row, $col)
        $link = [string]$cell.Text
        if ($link) {
            $code = $link.Substring($link.LastIndexOf("=") + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

Update-MatchCodeColumn $battingSheet 4 3
Update-MatchCodeColumn $bowlingSheet 2 3

Write-Host "Sheets now: $(($wb.Worksheets | ForEach-Object { $_.Name }) -join ', ')"
